$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").ClearContents()
$ws.Range("A1").Value = 'questions = [
    {
        "title": "You are a marketing manager for an online retail company. You are looking to expand your email marketing list to send promotional materials and updates. What is the most appropriate method to accomplish this?",
        "ques_type": 2,
        "options": [
            "Include an email field on the computer sign-up list available at the register of a store.",
            "Use a customer\u2019s email from a recent purchase on your e-commerce site.",
            "Ask for a customer\u2019s email in a pop-up window on your site and offer a 10% discount.",
            "Ask for a reader\u2019s email to subscribe to your blog."
        ],
        "score": "Use a customer\u2019s email from a recent purchase on your e-commerce site."
    },
    {
        "title": "Your company\u2019s overall goal for the month is to increase sales of Product A. What is an appropriate goal to set for your email campaign for the month?",
        "ques_type": 2,
        "options": [
            "Get a 10% click-through rate on a promotional campaign linking to Product A\u2019s product page.",
            "Increase the email open rate by 5% for campaigns focused on promoting Product A compared to last month\u2019s campaign.",
            "Get a 2% conversion rate from a campaign promoting Product A.",
            "Gain 50 new email subscribers from Product A\u2019s page."
        ],
        "score": "Get a 2% conversion rate from a campaign promoting Product A."
    },
    {
        "title": "You are an email marketing manager for a B2B services company. Analyzing the statistics of the last months, you see that the unsubscribe rate is getting higher with every email sent. You decide to run an A/B test to improve this metric.What would you test first?",
        "ques_type": 2,
        "options": [
            "Offer various discount offers and pricing strategies.",
            "Alter the design and the layout of your emails. ",
            "Change the image''s style and color in the emails. ",
            "Provide an option to adjust email frequency on the Unsubscribe page."
        ],
        "score": "Provide an option to adjust email frequency on the Unsubscribe page."
    },
    {
        "title": "You are a digital marketing specialist for a tech company. Your recent email campaign has an open rate of 12% and a click-through rate of 0.3%. You are tasked with identifying potential reasons for the low click-through rate. What might be causing this issue?",
        "ques_type": 15,
        "options": [
            "A broken link within the email",
            "The absence of a link from the first third of the email body",
            "An unclear call to action ",
            "Unengaging content",
            "The absence of any clickable links in the email"
        ],
        "score": [
            "The absence of a link from the first third of the email body",
            "An unclear call to action",
            "Unengaging content"
        ]
    }
]'
$ws.Range("A1").Style = "Normal"
$ws.Rows(1).AutoFit()
